$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("D")
$ws2 = $wb.Worksheets.Item("p-value")

$ws1.Range("A10:A11").UnMerge() | Out-Null
$ws1.Range("A12:A28").UnMerge() | Out-Null
$ws1.Range("A8:A31").ClearContents() | Out-Null
$ws1.Cells.Item(4, 3).Value = 0.06225618820524588
$ws1.Cells.Item(4, 4).Value = 0.05015661212735886
$ws1.Cells.Item(5, 3).Value = 0.1402560086200876
$ws1.Cells.Item(5, 4).Value = 0.05532542330078897
$ws1.Cells.Item(6, 3).Value = 0.05565144516167653
$ws1.Cells.Item(6, 4).Value = 0.03415325586226279
$ws1.Cells.Item(7, 3).Value = 0.092486356516447
$ws1.Cells.Item(7, 4).Value = 0.05431935839941999
$ws1.Cells.Item(8, 2).Value = "ASM1 mu H [/d]"
$ws1.Cells.Item(8, 3).Value = 0.113268350111243
$ws1.Cells.Item(8, 4).Value = 0.07991902603142018
$ws1.Cells.Item(9, 2).Value = "ASM1 K S [g COD/m3]"
$ws1.Cells.Item(9, 3).Value = 0.06982869571290319
$ws1.Cells.Item(9, 4).Value = 0.0476068282072901
$ws1.Cells.Item(10, 2).Value = "ASM1 K O H [g O2/m3]"
$ws1.Cells.Item(10, 3).Value = 0.2456226117668187
$ws1.Cells.Item(10, 4).Value = 0.09022406328026035
$ws1.Cells.Item(11, 2).Value = "ASM1 K NO [g N/m3]"
$ws1.Cells.Item(11, 3).Value = 0.1170995001546428
$ws1.Cells.Item(11, 4).Value = 0.0829942446570622
$ws1.Cells.Item(12, 2).Value = "ASM1 b H [/d]"
$ws1.Cells.Item(12, 3).Value = 0.1263980205724776
$ws1.Cells.Item(12, 4).Value = 0.05125635918553548
$ws1.Cells.Item(13, 2).Value = "ASM1 mu A [/d]"
$ws1.Cells.Item(13, 3).Value = 0.1731300695393641
$ws1.Cells.Item(13, 4).Value = 0.1445963724639629
$ws1.Cells.Item(14, 2).Value = "ASM1 K NH [g N/m3]"
$ws1.Cells.Item(14, 3).Value = 0.1455837016491904
$ws1.Cells.Item(14, 4).Value = 0.1610885051993597
$ws1.Cells.Item(15, 2).Value = "ASM1 K O A [g COD/m3]"
$ws1.Cells.Item(15, 3).Value = 0.1684508784707326
$ws1.Cells.Item(15, 4).Value = 0.1135631397371197
$ws1.Cells.Item(16, 2).Value = "ASM1 b A [/d]"
$ws1.Cells.Item(16, 3).Value = 0.1248216619608704
$ws1.Cells.Item(16, 4).Value = 0.1101661432685297
$ws1.Cells.Item(17, 2).Value = "ASM1 eta g"
$ws1.Cells.Item(17, 3).Value = 0.1415131047280781
$ws1.Cells.Item(17, 4).Value = 0.0811531866189295
$ws1.Cells.Item(18, 2).Value = "ASM1 k a [m3/g COD/d]"
$ws1.Cells.Item(18, 3).Value = 0.1354770480190759
$ws1.Cells.Item(18, 4).Value = 0.08178044975581542
$ws1.Cells.Item(19, 2).Value = "ASM1 k h [g X_S/g X_BH COD/d]"
$ws1.Cells.Item(19, 3).Value = 0.1546527521425507
$ws1.Cells.Item(19, 4).Value = 0.05821327761281572
$ws1.Cells.Item(20, 2).Value = "ASM1 K X [g X_S/g X_BH COD]"
$ws1.Cells.Item(20, 3).Value = 0.06406201674132753
$ws1.Cells.Item(20, 4).Value = 0.04058473958396976
$ws1.Cells.Item(21, 2).Value = "ASM1 eta h"
$ws1.Cells.Item(21, 3).Value = 0.08142191537548263
$ws1.Cells.Item(21, 4).Value = 0.05716240820166917
$ws1.Cells.Item(22, 2).Value = "ASM1 Y H [g COD/g COD]"
$ws1.Cells.Item(22, 3).Value = 0.076652931727709
$ws1.Cells.Item(22, 4).Value = 0.04245023644561751
$ws1.Cells.Item(23, 2).Value = "ASM1 Y A [g COD/g N]"
$ws1.Cells.Item(23, 3).Value = 0.08604124472468598
$ws1.Cells.Item(23, 4).Value = 0.05822142388732073
$ws1.Cells.Item(24, 2).Value = "ASM1 f pobs"
$ws1.Cells.Item(24, 3).Value = 0.07346030669154253
$ws1.Cells.Item(24, 4).Value = 0.1019302597439626
$ws1.Cells.Item(25, 2).Value = "Aerobic zone hydraulic retention time [hr]"
$ws1.Cells.Item(25, 1).Value = "CSTR-O1"
$ws1.Cells.Item(25, 3).Value = 0.2506410192455428
$ws1.Cells.Item(25, 4).Value = 0.2291261898652199
$ws1.Cells.Item(26, 2).Value = "O1 and o2 kla"
$ws1.Cells.Item(26, 3).Value = 0.5121369636140516
$ws1.Cells.Item(26, 4).Value = 0.5030813283315207
$ws1.Cells.Item(27, 2).Value = "Saturation DO [mg/L]"
$ws1.Cells.Item(27, 3).Value = 0.2627630174297373
$ws1.Cells.Item(27, 4).Value = 0.2688963019986885
$ws1.Cells.Item(28, 2).Value = "Internal recirculation rate as a fraction of influent"
$ws1.Cells.Item(28, 1).Value = "CSTR-O3"
$ws1.Cells.Item(28, 3).Value = 0.06799293631710748
$ws1.Cells.Item(28, 4).Value = 0.05458003918358037
$ws1.Cells.Item(29, 2).Value = "O3 kla"
$ws1.Cells.Item(29, 3).Value = 0.06570821402560086
$ws1.Cells.Item(29, 4).Value = 0.08492491171475006
$ws1.Cells.Item(30, 2).Value = "Sludge recycling as a fraction of influent"
$ws1.Cells.Item(30, 1).Value = "Flat bottom circular clarifier-C1"
$ws1.Cells.Item(30, 3).Value = 0.1079805648950923
$ws1.Cells.Item(30, 4).Value = 0.1674833306857941
$ws1.Cells.Item(31, 2).Value = "Waste sludge flowrate [m3/d]"
$ws1.Cells.Item(31, 3).Value = 0.12956071474893
$ws1.Cells.Item(31, 4).Value = 0.1845742145973093
$ws1.Range("A7:A24").Merge() | Out-Null
$ws1.Range("A25:A27").Merge() | Out-Null
$ws1.Range("A28:A29").Merge() | Out-Null
$ws1.Range("A30:A31").Merge() | Out-Null

$ws2.Range("A10:A11").UnMerge() | Out-Null
$ws2.Range("A12:A28").UnMerge() | Out-Null
$ws2.Range("A8:A31").ClearContents() | Out-Null
$ws2.Cells.Item(4, 3).Value = 0.8062248012750971
$ws2.Cells.Item(4, 4).Value = 0.5451617631770604
$ws2.Cells.Item(5, 3).Value = 0.03473006018611813
$ws2.Cells.Item(5, 4).Value = 0.4203241132540868
$ws2.Cells.Item(6, 3).Value = 0.8961722686892706
$ws2.Cells.Item(6, 4).Value = 0.9243211649889043
$ws2.Cells.Item(7, 3).Value = 0.3349649328115392
$ws2.Cells.Item(7, 4).Value = 0.4433230517368963
$ws2.Cells.Item(8, 2).Value = "ASM1 mu H [/d]"
$ws2.Cells.Item(8, 3).Value = 0.1404084612499558
$ws2.Cells.Item(8, 4).Value = 0.08109185701253037
$ws2.Cells.Item(9, 2).Value = "ASM1 K S [g COD/m3]"
$ws2.Cells.Item(9, 3).Value = 0.683263662940607
$ws2.Cells.Item(9, 4).Value = 0.611270278843504
$ws2.Cells.Item(10, 2).Value = "ASM1 K O H [g O2/m3]"
$ws2.Cells.Item(10, 3).Value = [double]"8.256633459367e-06"
$ws2.Cells.Item(10, 4).Value = 0.03390437689712211
$ws2.Cells.Item(11, 2).Value = "ASM1 K NO [g N/m3]"
$ws2.Cells.Item(11, 3).Value = 0.1172684402856004
$ws2.Cells.Item(11, 4).Value = 0.06320515193315844
$ws2.Cells.Item(12, 2).Value = "ASM1 b H [/d]"
$ws2.Cells.Item(12, 3).Value = 0.07386994212750495
$ws2.Cells.Item(12, 4).Value = 0.5172952279790757
$ws2.Cells.Item(13, 2).Value = "ASM1 mu A [/d]"
$ws2.Cells.Item(13, 3).Value = 0.004230074828731686
$ws2.Cells.Item(13, 4).Value = [double]"5.983344002405235e-05"
$ws2.Cells.Item(14, 2).Value = "ASM1 K NH [g N/m3]"
$ws2.Cells.Item(14, 3).Value = 0.02541676946673458
$ws2.Cells.Item(14, 4).Value = [double]"4.878938813632685e-06"
$ws2.Cells.Item(15, 2).Value = "ASM1 K O A [g COD/m3]"
$ws2.Cells.Item(15, 3).Value = 0.005866027024870765
$ws2.Cells.Item(15, 4).Value = 0.003200086402620198
$ws2.Cells.Item(16, 2).Value = "ASM1 b A [/d]"
$ws2.Cells.Item(16, 3).Value = 0.08008577472594217
$ws2.Cells.Item(16, 4).Value = 0.004663045448409871
$ws2.Cells.Item(17, 2).Value = "ASM1 eta g"
$ws2.Cells.Item(17, 3).Value = 0.03229832419904825
$ws2.Cells.Item(17, 4).Value = 0.07340356981168983
$ws2.Cells.Item(18, 2).Value = "ASM1 k a [m3/g COD/d]"
$ws2.Cells.Item(18, 3).Value = 0.04545476100032474
$ws2.Cells.Item(18, 4).Value = 0.06978419152192343
$ws2.Cells.Item(19, 2).Value = "ASM1 k h [g X_S/g X_BH COD/d]"
$ws2.Cells.Item(19, 3).Value = 0.01457650018623719
$ws2.Cells.Item(19, 4).Value = 0.3583380859833175
$ws2.Cells.Item(20, 2).Value = "ASM1 K X [g X_S/g X_BH COD]"
$ws2.Cells.Item(20, 3).Value = 0.7781994259899252
$ws2.Cells.Item(20, 4).Value = 0.7930410942192309
$ws2.Cells.Item(21, 2).Value = "ASM1 eta h"
$ws2.Cells.Item(21, 3).Value = 0.4916495304966206
$ws2.Cells.Item(21, 4).Value = 0.3802261827724651
$ws2.Cells.Item(22, 2).Value = "ASM1 Y H [g COD/g COD]"
$ws2.Cells.Item(22, 3).Value = 0.5684565624176354
$ws2.Cells.Item(22, 4).Value = 0.7464416338162054
$ws2.Cells.Item(23, 2).Value = "ASM1 Y A [g COD/g N]"
$ws2.Cells.Item(23, 3).Value = 0.4220608792171776
$ws2.Cells.Item(23, 4).Value = 0.3580958582848331
$ws2.Cells.Item(24, 2).Value = "ASM1 f pobs"
$ws2.Cells.Item(24, 3).Value = 0.6220596507502912
$ws2.Cells.Item(24, 4).Value = 0.01110895792203313
$ws2.Cells.Item(25, 2).Value = "Aerobic zone hydraulic retention time [hr]"
$ws2.Cells.Item(25, 1).Value = "CSTR-O1"
$ws2.Cells.Item(25, 3).Value = [double]"4.932821914471342e-06"
$ws2.Cells.Item(25, 4).Value = [double]"8.509111524568462e-12"
$ws2.Cells.Item(26, 2).Value = "O1 and o2 kla"
$ws2.Cells.Item(26, 3).Value = [double]"9.360553488726321e-25"
$ws2.Cells.Item(26, 4).Value = [double]"4.235134478122423e-57"
$ws2.Cells.Item(27, 2).Value = "Saturation DO [mg/L]"
$ws2.Cells.Item(27, 3).Value = [double]"1.358552934530584e-06"
$ws2.Cells.Item(27, 4).Value = [double]"3.989584221254844e-16"
$ws2.Cells.Item(28, 2).Value = "Internal recirculation rate as a fraction of influent"
$ws2.Cells.Item(28, 1).Value = "CSTR-O3"
$ws2.Cells.Item(28, 3).Value = 0.7141986555968359
$ws2.Cells.Item(28, 4).Value = 0.4372178424932672
$ws2.Cells.Item(29, 2).Value = "O3 kla"
$ws2.Cells.Item(29, 3).Value = 0.7517594452307359
$ws2.Cells.Item(29, 4).Value = 0.05380778794020067
$ws2.Cells.Item(30, 2).Value = "Sludge recycling as a fraction of influent"
$ws2.Cells.Item(30, 1).Value = "Flat bottom circular clarifier-C1"
$ws2.Cells.Item(30, 3).Value = 0.1781934317256087
$ws2.Cells.Item(30, 4).Value = [double]"1.712805685741369e-06"
$ws2.Cells.Item(31, 2).Value = "Waste sludge flowrate [m3/d]"
$ws2.Cells.Item(31, 3).Value = 0.06260109546839047
$ws2.Cells.Item(31, 4).Value = [double]"8.543660839223266e-08"
$ws2.Range("A7:A24").Merge() | Out-Null
$ws2.Range("A25:A27").Merge() | Out-Null
$ws2.Range("A28:A29").Merge() | Out-Null
$ws2.Range("A30:A31").Merge() | Out-Null
